# Automated daily market-data refresh: append the newest trading day to each sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Prix Spot" sheet: add a new day column BY (29-aug) with its 25 hourly rows.
# ---------------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Copy the formatting of the last existing day column (BX) into the new one
# (BY) first, so the header keeps the same bold/border/centered style and the
# body cells keep the plain numeric formatting used throughout the sheet.
$wsPrix.Range("BX1:BX25").Copy($wsPrix.Range("BY1:BY25"))

$wsPrix.Range("BY1").Value = "29-aug"
$wsPrix.Range("BY2").Value = 82.7
$wsPrix.Range("BY3").Value = 51.1
$wsPrix.Range("BY4").Value = 51.93
$wsPrix.Range("BY5").Value = 25.04
$wsPrix.Range("BY6").Value = 17.87
$wsPrix.Range("BY7").Value = 21.33
$wsPrix.Range("BY8").Value = 30.57
$wsPrix.Range("BY9").Value = 45.47
$wsPrix.Range("BY10").Value = 57.53
$wsPrix.Range("BY11").Value = 42.47
$wsPrix.Range("BY12").Value = 5
$wsPrix.Range("BY13").Value = 0.83
$wsPrix.Range("BY14").Value = 3.08
$wsPrix.Range("BY15").Value = 0.05
$wsPrix.Range("BY16").Value = 0
$wsPrix.Range("BY17").Value = 0.62
$wsPrix.Range("BY18").Value = 5.66
$wsPrix.Range("BY19").Value = 3.42
$wsPrix.Range("BY20").Value = 22.36
$wsPrix.Range("BY21").Value = 53.55
$wsPrix.Range("BY22").Value = 81.8
$wsPrix.Range("BY23").Value = 63.5
$wsPrix.Range("BY24").Value = 94.39
$wsPrix.Range("BY25").Value = 86

# ---------------------------------------------------------------------------
# "Gaz" sheet: append row 74 (2025-08-27, 31.4)
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

# Force the date column to plain text first so the new date string isn't
# silently converted into a date serial number (matching how every other
# row in the column already stores its date as literal text), then drop the
# explicit style again so the cell matches the unstyled body cells above it.
$wsGaz.Range("A74").NumberFormat = "@"
$wsGaz.Range("A74").Value = "2025-08-27"
$wsGaz.Range("A74").Style = "Normal"
$wsGaz.Range("B74").Value = 31.4

# ---------------------------------------------------------------------------
# "CO2" sheet: append row 74 (2025-08-27, 72.36)
# ---------------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")

$wsCo2.Range("A74").NumberFormat = "@"
$wsCo2.Range("A74").Value = "2025-08-27"
$wsCo2.Range("A74").Style = "Normal"
$wsCo2.Range("B74").Value = 72.36
